# Update column G ("K") values for rows 2-25 in the active worksheet.
# These values represent a recalculated statistic (K, replacing the old
# "Strike#" values) for each saved record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 4
    4  = 2
    5  = 3
    6  = 1
    7  = 2
    8  = 8
    9  = 2
    10 = 6
    11 = 5
    12 = 4
    13 = 3
    14 = 3
    15 = 3
    16 = 3
    17 = 1
    18 = 3
    19 = 2
    20 = 1
    21 = 2
    22 = 3
    23 = 2
    24 = 3
    25 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
